# Informe 1 casi completo
# Continue the "Regla Falsa" (False Position) root-finding table:
# update the initial interval (x_i, x_s) used to search for the root of
# f(x) = x^2 - 2, and extend the iteration table from 7 iterations (rows
# 2-8) up to 15 iterations (rows 2-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function F($x) {
    return ($x * $x) - 2
}

# New initial bracket for the root search.
$xi = -4.0
$xs = 1.0
$fi = F($xi)
$fs = F($xs)

for ($n = 1; $n -le 15; $n++) {
    $row = $n + 1

    $xm = $xs - ($fs * ($xi - $xs) / ($fi - $fs))
    $fm = F($xm)

    if ($n -eq 1) {
        $err = 1.0005
    } else {
        $err = [Math]::Abs($xm - $xs)
    }

    $ws.Cells.Item($row, 1).Value = $n
    $ws.Cells.Item($row, 2).Value = $xm
    $ws.Cells.Item($row, 3).Value = $xi
    $ws.Cells.Item($row, 4).Value = $xs
    $ws.Cells.Item($row, 5).Value = $fm
    $ws.Cells.Item($row, 6).Value = $fi
    $ws.Cells.Item($row, 7).Value = $fs
    $ws.Cells.Item($row, 8).Value = $err

    if (($fm * $fs) -gt 0) {
        $xs = $xm
        $fs = $fm
    } else {
        $xi = $xm
        $fi = $fm
    }
}
